# task3 and task4 - re-run after loss function change
# (this workbook only contains the "task4" sheet/output)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Re-run results: updated test_error (column D) values ---
$ws.Range("D2").Value  = 0.10390000000000001
$ws.Range("D3").Value  = 0.10929999999999999
$ws.Range("D4").Value  = 0.28899999999999998
$ws.Range("D5").Value  = 0.1099
$ws.Range("D6").Value  = 0.22900000000000001
$ws.Range("D7").Value  = 0.8266
$ws.Range("D8").Value  = 0.1174
$ws.Range("D9").Value  = 0.073899999999999993
$ws.Range("D10").Value = 0.20880000000000001
$ws.Range("D11").Value = 0.068599999999999994
$ws.Range("D12").Value = 0.1555
$ws.Range("D13").Value = 0.4456

# Row 11 had been highlighted (best result); the highlight moved off of
# A11:C11 on the re-run, leaving only D11 styled.
$ws.Range("A11:C11").Style = "Normal"

# Selection moved to F16 before the file was saved
$ws.Range("F16").Select()

# Rename the sheet (and, with it, the filter-database defined name that
# references it) to "task4 output"
$ws.Name = "task4 output"
